$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook calc settings: concurrentCalc="0" (disable multi-threaded calc) ---
$excel.MultiThreadedCalculation.Enabled = $false

# --- Restructure the table: insert 3 new blank columns before the old column B ---
# This shifts the old B..AK content (and its per-column styles) to E..AN,
# leaving A untouched for now.
$ws.Range("B1:D1").EntireColumn.Insert()

# Move the "#" header (old A2) up to the new header cell A1
$ws.Range("A1").Value = $ws.Range("A2").Value2
$ws.Range("A2").ClearContents()

# --- New leading "test" columns (A:D) on the field-header rows 2-4 ---
$ws.Range("A2").Value = "測試用"
$ws.Range("B2").Value = "TstByte"
$ws.Range("C2").Value = "TestUINT"
$ws.Range("D2").Value = "testnon"

$ws.Range("A3").Value = "string"
$ws.Range("B3").Value = "byte"
$ws.Range("C3").Value = "uint"
$ws.Range("D3").Value = "string"

$ws.Range("A4").Value = "C"
$ws.Range("B4").Value = "C"
$ws.Range("C4").Value = "C"
$ws.Range("D4").Value = "N"

# Last header cell changes from "End" to "EOC"
$ws.Range("AN2").Value = "EOC"

# --- Clear out the old sample-data rows (5-8); rows 5-7 are rebuilt below ---
$ws.Range("A5:AN8").ClearContents()

# --- Row 5: full sample row ---
$ws.Range("A5").Value = 123
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = "0aa"
$ws.Range("E5").Value = 10
$ws.Range("G5").Value = 14
$ws.Range("H5").Value = 16
$ws.Range("I5").Value = 18
$ws.Range("J5").Value = 20
$ws.Range("K5").Value = 22
$ws.Range("L5").Value = 24
$ws.Range("M5").Value = 26
$ws.Range("N5").Value = 28
$ws.Range("O5").Value = 30
$ws.Range("P5").Value = 32
$ws.Range("Q5").Value = 34
$ws.Range("R5").Value = 36
$ws.Range("S5").Value = 38
$ws.Range("T5").Value = 40
$ws.Range("U5").Value = 42
$ws.Range("V5").Value = 44
$ws.Range("W5").Value = 46
$ws.Range("X5").Value = 48
$ws.Range("Y5").Value = 50
$ws.Range("Z5").Value = 52
$ws.Range("AA5").Value = 54
$ws.Range("AB5").Value = 56
$ws.Range("AC5").Value = 58
$ws.Range("AD5").Value = 60
$ws.Range("AE5").Value = 62
$ws.Range("AF5").Value = 64
$ws.Range("AG5").Value = 66
$ws.Range("AH5").Value = 68
$ws.Range("AI5").Value = 70
$ws.Range("AJ5").Value = 72
$ws.Range("AK5").Value = 74
$ws.Range("AL5").Value = 76
$ws.Range("AM5").Value = 78

# --- Row 6: sparse sample row ---
$ws.Range("A6").Value = "a23"
$ws.Range("B6").Value = 1
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 9
$ws.Range("I6").Value = 123
$ws.Range("L6").Value = 50
$ws.Range("O6").Value = 98
$ws.Range("R6").Value = 31
$ws.Range("U6").Value = 56465
$ws.Range("X6").Value = 12
$ws.Range("AA6").Value = 50
$ws.Range("AD6").Value = 1254
$ws.Range("AJ6").Value = 66
$ws.Range("AM6").Value = 77

# --- Row 7: terminator marker ---
$ws.Range("A7").Value = "EOR"

# --- Sheet view: restore selection/scroll position ---
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 16
$null = $ws.Range("AL6").Select()

# --- Page setup (paper size / orientation) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
